$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.960.01"
$ws.Range("E2").Value = "  -1.02%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.215.91"
$ws.Range("E3").Value = "  -1.67%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.48"
$ws.Range("E5").Value = "  -2.25%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.627"
$ws.Range("E6").Value = "  -0.77%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "73.26"
$ws.Range("E7").Value = "  -3.08%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.606"
$ws.Range("E9").Value = "  -2.84%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.64"
$ws.Range("E10").Value = "  -3.20%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0957"
$ws.Range("E11").Value = "  +0.62%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.02"
$ws.Range("E12").Value = "  -2.89%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.104"
$ws.Range("E13").Value = "  +0.72%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.547.48"
$ws.Range("E14").Value = "  -1.57%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.25"
$ws.Range("E15").Value = "  -2.50%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.836"
$ws.Range("E16").Value = "  -2.55%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.217.53"
$ws.Range("E17").Value = "  -1.90%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.917.35"
$ws.Range("E18").Value = "  -0.76%  "

$ws.Range("E19").Value = "  +7.07%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.18"
$ws.Range("E20").Value = "  -0.03%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.58"
$ws.Range("E21").Value = "  +0.45%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.92"
$ws.Range("E22").Value = "  +17.90%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "228.99"
$ws.Range("E23").Value = "  -1.41%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.07"
$ws.Range("E24").Value = "  -7.53%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.76"
$ws.Range("E25").Value = "  +1.95%  "

$ws.Range("E26").Value = "  -0.10%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.74"
$ws.Range("E27").Value = "  +3.58%  "

$ws.Range("E28").Value = "  -1.23%  "

$ws.Range("E29").Value = "  -0.99%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "167.74"
$ws.Range("E30").Value = "  -0.79%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.50"
$ws.Range("E31").Value = "  -1.16%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.68"
$ws.Range("E32").Value = "  +7.98%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0797"
$ws.Range("E33").Value = "  -3.41%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "30.38"
$ws.Range("E34").Value = "  -1.24%  "

$ws.Range("E35").Value = "  -0.08%  "

$ws.Range("E36").Value = "  -9.82%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.26"
$ws.Range("E37").Value = "  -3.84%  "

$ws.Range("E38").Value = "  -4.11%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.77"
$ws.Range("E39").Value = "  -1.50%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "65.76"
$ws.Range("E40").Value = "  +4.24%  "

$ws.Range("E41").Value = "  -3.43%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.64"
$ws.Range("E42").Value = "  -2.91%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "8.86"
$ws.Range("E43").Value = "  +0.64%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D44").Value = "0.198"
$ws.Range("E44").Value = "  -3.06%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "104.74"
$ws.Range("E45").Value = "  -3.79%  "

$ws.Range("E46").Value = "  -2.27%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.40"
$ws.Range("E47").Value = "  +3.29%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.11"
$ws.Range("E48").Value = "  -1.27%  "

$ws.Range("E49").Value = "  -1.89%  "

$ws.Range("E50").Value = "  -0.29%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.422.71"
$ws.Range("E51").Value = "  -1.84%  "
